# Sync attendance_reports: normalize the 'Recorded By' (column G) ordering.
# For each affected session row, the last author in the comma-separated
# 'Recorded By' list is moved to the front (the list is right-rotated by one),
# matching the upstream main-repo sync.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in column G whose 'Recorded By' value needs reordering.
$rowsToFix = @(
    2, 3, 6, 7, 10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21,
    22, 24, 26, 28, 29, 32, 33, 36, 37, 38, 39, 40, 41, 43, 44,
    45, 46, 47, 48, 50, 52, 54, 55, 58, 59, 62, 63, 64, 65, 66,
    67, 69, 70, 71, 72, 73, 74, 76, 78, 83, 84, 85, 86, 87, 90,
    92, 93, 94, 96, 99, 101, 109, 110, 111, 112, 113, 116, 118, 119, 120,
    122, 125, 127, 135, 136, 137, 138, 139, 142, 144, 145, 146, 148, 151, 153
)

$updatedCount = 0
foreach ($row in $rowsToFix) {
    $cell = $ws.Range("G$row")
    $current = [string]$cell.Text
    $parts = @($current -split ", ")
    if ($parts.Count -gt 1) {
        $lastAuthor = $parts[$parts.Count - 1]
        $rest = $parts[0..($parts.Count - 2)]
        $newValue = ($lastAuthor + ", " + ($rest -join ", "))
        $cell.Value = $newValue
        $updatedCount++
    }
}

Write-Output ("Updated " + $updatedCount + " cells in column G")
